# Insert a new column before column A, shifting existing data (B:J -> C:K)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").EntireColumn.Insert()

# New header for column A (copy formatting from the adjacent header cell)
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A1").Value = "Scrip_Symbol"

# New Scrip_Symbol values for rows 2-12
$values = @(
    "BSE BSE - MITHIL DEEPAK KOTWAL",
    "544223 CEIGALL - MITHIL DEEPAK KOTWAL",
    "544271 GARUDA - MITHIL DEEPAK KOTWAL",
    "500116 IDBI - MITHIL DEEPAK KOTWAL",
    "543398 LATENTVIEW - MITHIL DEEPAK KOTWAL",
    "532461 PNB - MITHIL DEEPAK KOTWAL",
    "532461 PNB - MITHIL DEEPAK KOTWAL",
    "532461 PNB - MITHIL DEEPAK KOTWAL",
    "544243 STYLEBAAZA - MITHIL DEEPAK KOTWAL",
    "521064 TRIDENT - MITHIL DEEPAK KOTWAL",
    "543238 UTIAMC - MITHIL DEEPAK KOTWAL"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
